$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# K2: was a literal number (40675536005), becomes the text value "40676340007".
# Writing directly via .Value keeps K2's existing style (no quotePrefix to lose).
$ws.Range("K2").Value = "40676340007"

# B2: text value changes from "1035655536" to "1035103510".
# B2's original style (quotePrefix text style) gets stripped if we assign
# .Value directly, so stage the new text in a scratch cell, copy it, and
# paste-values-only into B2 so its existing number format/border/quotePrefix
# style is preserved. Clean up the scratch cell afterwards.
$ws.Range("ZZ1").Value = "1035103510"
$ws.Range("ZZ1").Copy()
$ws.Range("B2").PasteSpecial(-4163)
$ws.Range("ZZ1").Clear()

# Update the active selection shown when the sheet is reopened.
$ws.Range("D20").Select()
